$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# Insert a new blank row at position 4, shifting existing rows 4-12 down to 5-13
$ws.Rows.Item(4).Insert()

# The sheet stores every value as text (even numeric-looking ones), so force
# the new cells to Text format before writing, otherwise Excel will
# auto-convert numeric-looking strings ("151", "300", ...) into real numbers.
$ws.Range("A4:O4").NumberFormat = "@"

# Fill in the new row 4 values
$ws.Cells.Item(4, 1).Value = "1"
$ws.Cells.Item(4, 2).Value = "151"
$ws.Cells.Item(4, 3).Value = "300"
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(4, 7).Value = ""
$ws.Cells.Item(4, 8).Value = ""
$ws.Cells.Item(4, 9).Value = ""
$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 11).Value = ""
$ws.Cells.Item(4, 12).Value = "11\8\2024"
$ws.Cells.Item(4, 13).Value = "4300 ساعه تشغيل"
$ws.Cells.Item(4, 14).Value = "تم التشحيم شحم Ep3"
$ws.Cells.Item(4, 15).Value = "حسام"

# The new row had no styling in the source workbook (no "s" attribute), so
# drop the temporary Text number format again now that the values are typed
# as text - this keeps the cells styled exactly like the rest of the sheet.
$ws.Range("A4:O4").ClearFormats()
